$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts Resolved on..Problem right by one)
$ws.Columns("D:D").Insert()

# New defect row D9 (order chosen to match shared-string insertion order)
$ws.Range("A23").Value = "D9"
$ws.Range("I23").Value = "On analysis, not more than 25 tweets are analysed and saved (even if they can be viewed in the filter)."
$ws.Range("D23").Value = "johannes"

# New column header
$ws.Range("D1").Value = "Assignee"

$ws.Range("C23").Value = "assigned"

# Selection / view
$ws.Range("C23").Select()
